# docs: add responsibility matrix
# Halve the "Carico di Lavoro (ore/uomo)" (column C) workload values for the
# "Backend" section (rows 72-95, skipping the blank separator rows), update
# the grand-total row's hour-columns (C168/D168) to plain number formatting
# (they represent hours, not currency) and refresh the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Halve column C workload values for rows 72-95 (Backend section) -------
$backendRows = @(72, 73, 74, 76, 77, 78, 80, 81, 82, 84, 85, 86, 88, 89, 90, 91, 93, 94, 95)
foreach ($r in $backendRows) {
    $cell = $ws.Cells.Item($r, 3)   # column C
    $cell.Value = $cell.Value() / 2
}

# --- Grand-total row: hours columns should show plain numbers, not currency
$ws.Range("C168").NumberFormat = "General"
$ws.Range("D168").NumberFormat = "General"

# --- Update the visible selection / scroll position -------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 147
$win.ScrollColumn = 1
$ws.Range("C177").Select()
